$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 287
    3  = 288
    4  = 290
    5  = 293
    6  = 295
    7  = 298
    8  = 300
    9  = 302
    10 = 303
    11 = 305
    12 = 308
    13 = 310
    14 = 311
    15 = 18
    16 = 41
    17 = 58
    18 = 124
    19 = 130
    20 = 169
    21 = 205
    22 = 213
    23 = 237
    24 = 372
    25 = 463
    26 = 478
    27 = 505
}

foreach ($row in $values.Keys) {
    $ws.Range("A$row").Value = $values[$row]
}
